$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
# Row 17
$ws.Cells.Item(17, 8).Value = 3647402.2  # H17: was 3414673.2
$ws.Cells.Item(17, 10).Value = 3647402.2  # J17: was 3414673.2
$ws.Cells.Item(17, 12).Value = 10942206.6  # L17: was 10244019.6
$ws.Cells.Item(17, 14).Value = -10942542.6  # N17: was -10244355.6

# Row 39
$ws.Cells.Item(39, 8).Value = 112  # H39: was 109
$ws.Cells.Item(39, 10).Value = 0  # J39: was 100
$ws.Cells.Item(39, 12).Value = 0  # L39: was 300
$ws.Cells.Item(39, 14).ClearContents()  # N39: was -892

# Row 80
$ws.Cells.Item(80, 8).Value = 3451.75  # H80: was 3452.875
$ws.Cells.Item(80, 9).Value = 368.41177  # I80: was 357.16666
$ws.Cells.Item(80, 10).Value = 6946.2  # J80: was 7433.0713
$ws.Cells.Item(80, 11).Value = 1105.23531  # K80: was 1071.49998
$ws.Cells.Item(80, 12).Value = 20838.6  # L80: was 22299.2139
$ws.Cells.Item(80, 13).Value = -107.23531  # M80: was -73.49998000000005
$ws.Cells.Item(80, 14).Value = -22834.6  # N80: was -24295.2139

# Row 83
$ws.Cells.Item(83, 8).Value = 3451.75  # H83: was 3452.875
$ws.Cells.Item(83, 9).Value = 368.41177  # I83: was 357.16666
$ws.Cells.Item(83, 10).Value = 6946.2  # J83: was 7433.0713
$ws.Cells.Item(83, 11).Value = 3315.70593  # K83: was 3214.49994
$ws.Cells.Item(83, 12).Value = 62515.8  # L83: was 66897.64169999999
$ws.Cells.Item(83, 13).Value = 1676.29407  # M83: was 1777.50006
$ws.Cells.Item(83, 14).Value = -72499.79999999999  # N83: was -76881.64169999999

# Row 125
$ws.Cells.Item(125, 8).Value = 3933.1428  # H125: was 6508
$ws.Cells.Item(125, 9).Value = 4422  # I125: was 8344
$ws.Cells.Item(125, 11).Value = 39798  # K125: was 75096
$ws.Cells.Item(125, 13).Value = -37338  # M125: was -72636

# Row 132
$ws.Cells.Item(132, 8).Value = 6538155  # H132: was 6668891.5
$ws.Cells.Item(132, 9).Value = 1489.375  # I132: was 1474.025
$ws.Cells.Item(132, 10).Value = 30307848  # J132: was 33338562
$ws.Cells.Item(132, 11).Value = 4468.125  # K132: was 4422.075000000001
$ws.Cells.Item(132, 12).Value = 90923544  # L132: was 100015686
$ws.Cells.Item(132, 13).Value = -1938.125  # M132: was -1892.075000000001
$ws.Cells.Item(132, 14).Value = -90928604  # N132: was -100020746

# Row 133
$ws.Cells.Item(133, 8).Value = 52580  # H133: was 50472
$ws.Cells.Item(133, 10).Value = 52580  # J133: was 50472
$ws.Cells.Item(133, 12).Value = 52580  # L133: was 50472
$ws.Cells.Item(133, 14).Value = -62700  # N133: was -60592

# Row 134
$ws.Cells.Item(134, 8).Value = 56888.832  # H134: was 57909.4
$ws.Cells.Item(134, 10).Value = 56888.832  # J134: was 57909.4
$ws.Cells.Item(134, 12).Value = 56888.832  # L134: was 57909.4
$ws.Cells.Item(134, 14).Value = -67028.83199999999  # N134: was -68049.39999999999

# Row 138
$ws.Cells.Item(138, 8).Value = 2652.2632  # H138: was 2594.5442
$ws.Cells.Item(138, 9).Value = 1494.258  # I138: was 1422.9117
$ws.Cells.Item(138, 10).Value = 4032.9614  # J138: was 3766.1765
$ws.Cells.Item(138, 11).Value = 4482.774  # K138: was 4268.7351
$ws.Cells.Item(138, 12).Value = 12098.8842  # L138: was 11298.5295
$ws.Cells.Item(138, 13).Value = 657.2259999999997  # M138: was 871.2649000000001
$ws.Cells.Item(138, 14).Value = -22378.8842  # N138: was -21578.5295

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
# Row 135
$ws.Cells.Item(135, 8).Value = 50804.75  # H135: was 53423.363
$ws.Cells.Item(135, 10).Value = 50804.75  # J135: was 53423.363
$ws.Cells.Item(135, 12).Value = 50804.75  # L135: was 53423.363
$ws.Cells.Item(135, 14).Value = -60944.75  # N135: was -63563.363

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
# Row 86
$ws.Cells.Item(86, 8).Value = 1970  # H86: was 2085.7144
$ws.Cells.Item(86, 9).Value = 1740  # I86: was 1800
$ws.Cells.Item(86, 10).Value = 2200  # J86: was 2133.3333
$ws.Cells.Item(86, 11).Value = 1740  # K86: was 1800
$ws.Cells.Item(86, 12).Value = 2200  # L86: was 2133.3333
$ws.Cells.Item(86, 13).Value = -617  # M86: was -677
$ws.Cells.Item(86, 14).Value = -4446  # N86: was -4379.3333

# Row 89
$ws.Cells.Item(89, 8).Value = 1970  # H89: was 2085.7144
$ws.Cells.Item(89, 9).Value = 1740  # I89: was 1800
$ws.Cells.Item(89, 10).Value = 2200  # J89: was 2133.3333
$ws.Cells.Item(89, 11).Value = 8700  # K89: was 9000
$ws.Cells.Item(89, 12).Value = 11000  # L89: was 10666.6665
$ws.Cells.Item(89, 13).Value = -3084  # M89: was -3384
$ws.Cells.Item(89, 14).Value = -22232  # N89: was -21898.6665

# Row 105
$ws.Cells.Item(105, 8).Value = 1706.9  # H105: was 1743.2222
$ws.Cells.Item(105, 9).Value = 1653.8  # I105: was 1722.25
$ws.Cells.Item(105, 11).Value = 1653.8  # K105: was 1722.25
$ws.Cells.Item(105, 13).Value = 93.20000000000005  # M105: was 24.75

# Row 137
$ws.Cells.Item(137, 8).Value = 51337.5  # H137: was 48608.89
$ws.Cells.Item(137, 10).Value = 42957.145  # J137: was 40935
$ws.Cells.Item(137, 12).Value = 42957.145  # L137: was 40935
$ws.Cells.Item(137, 14).Value = -53157.145  # N137: was -51135

# Row 140
$ws.Cells.Item(140, 8).Value = 47411.43  # H140: was 51711.668
$ws.Cells.Item(140, 10).Value = 47411.43  # J140: was 51711.668
$ws.Cells.Item(140, 12).Value = 47411.43  # L140: was 51711.668
$ws.Cells.Item(140, 14).Value = -57771.43  # N140: was -62071.668

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
# Row 31
$ws.Cells.Item(31, 8).Value = 9208998  # H31: was 9917766
$ws.Cells.Item(31, 9).Value = 1682.76  # I31: was 1808.2667
$ws.Cells.Item(31, 10).Value = 16634252  # J31: was 13937749
$ws.Cells.Item(31, 11).Value = 1682.76  # K31: was 1808.2667
$ws.Cells.Item(31, 12).Value = 16634252  # L31: was 13937749
$ws.Cells.Item(31, 13).Value = -1387.76  # M31: was -1513.2667
$ws.Cells.Item(31, 14).Value = -16634842  # N31: was -13938339

# Row 34
$ws.Cells.Item(34, 8).Value = 9208998  # H34: was 9917766
$ws.Cells.Item(34, 9).Value = 1682.76  # I34: was 1808.2667
$ws.Cells.Item(34, 10).Value = 16634252  # J34: was 13937749
$ws.Cells.Item(34, 11).Value = 1682.76  # K34: was 1808.2667
$ws.Cells.Item(34, 12).Value = 16634252  # L34: was 13937749
$ws.Cells.Item(34, 13).Value = -1480.76  # M34: was -1606.2667
$ws.Cells.Item(34, 14).Value = -16634656  # N34: was -13938153

# Row 132
$ws.Cells.Item(132, 8).Value = 5557818  # H132: was 8336256.5
$ws.Cells.Item(132, 9).Value = 8001469.5  # I132: was 14287639
$ws.Cells.Item(132, 10).Value = 4065.6365  # J132: was 4320.8
$ws.Cells.Item(132, 11).Value = 24004408.5  # K132: was 42862917
$ws.Cells.Item(132, 12).Value = 12196.9095  # L132: was 12962.4
$ws.Cells.Item(132, 13).Value = -24001878.5  # M132: was -42860387
$ws.Cells.Item(132, 14).Value = -17256.9095  # N132: was -18022.4

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
# Row 131
$ws.Cells.Item(131, 8).Value = 2440371.8  # H131: was 1696056
$ws.Cells.Item(131, 9).Value = 33333772  # I131: was 7692737.5
$ws.Cells.Item(131, 10).Value = 1419.0264  # J131: was 1341.7391
$ws.Cells.Item(131, 11).Value = 100001316  # K131: was 23078212.5
$ws.Cells.Item(131, 12).Value = 4257.0792  # L131: was 4025.2173
$ws.Cells.Item(131, 13).Value = -99996276  # M131: was -23073172.5
$ws.Cells.Item(131, 14).Value = -14337.0792  # N131: was -14105.2173

# Row 132
$ws.Cells.Item(132, 8).Value = 10895938  # H132: was 8716983
$ws.Cells.Item(132, 9).Value = 5004  # I132: was 2268
$ws.Cells.Item(132, 10).Value = 11886022  # J132: was 10895662
$ws.Cells.Item(132, 11).Value = 45036  # K132: was 20412
$ws.Cells.Item(132, 12).Value = 106974198  # L132: was 98060958
$ws.Cells.Item(132, 13).Value = -42506  # M132: was -17882
$ws.Cells.Item(132, 14).Value = -106979258  # N132: was -98066018

# Row 134
$ws.Cells.Item(134, 8).Value = 6725.1377  # H134: was 8169.0454
$ws.Cells.Item(134, 9).Value = 5914.304  # I134: was 7748.1763
$ws.Cells.Item(134, 10).Value = 9833.333000000001  # J134: was 9600
$ws.Cells.Item(134, 11).Value = 17742.912  # K134: was 23244.5289
$ws.Cells.Item(134, 12).Value = 29499.999  # L134: was 28800
$ws.Cells.Item(134, 13).Value = -12672.912  # M134: was -18174.5289
$ws.Cells.Item(134, 14).Value = -39639.999  # N134: was -38940

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
# Row 122
$ws.Cells.Item(122, 8).Value = 136483890  # H122: was 113737030
$ws.Cells.Item(122, 9).Value = 212964160  # I122: was 177470560
$ws.Cells.Item(122, 10).Value = 60003600  # J122: was 50003500
$ws.Cells.Item(122, 11).Value = 638892480  # K122: was 532411680
$ws.Cells.Item(122, 12).Value = 180010800  # L122: was 150010500
$ws.Cells.Item(122, 13).Value = -638890030  # M122: was -532409230
$ws.Cells.Item(122, 14).Value = -180015700  # N122: was -150015400

# Row 138
$ws.Cells.Item(138, 8).Value = 50429  # H138: was 39491.5
$ws.Cells.Item(138, 10).Value = 50429  # J138: was 39491.5
$ws.Cells.Item(138, 12).Value = 50429  # L138: was 39491.5
$ws.Cells.Item(138, 14).Value = -60709  # N138: was -49771.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
# Row 132
$ws.Cells.Item(132, 8).Value = 20007720  # H132: was 35726044
$ws.Cells.Item(132, 9).Value = 27787112  # I132: was 41679550
$ws.Cells.Item(132, 10).Value = 3569.7144  # J132: was 5002
$ws.Cells.Item(132, 11).Value = 83361336  # K132: was 125038650
$ws.Cells.Item(132, 12).Value = 10709.1432  # L132: was 15006
$ws.Cells.Item(132, 13).Value = -83358806  # M132: was -125036120
$ws.Cells.Item(132, 14).Value = -15769.1432  # N132: was -20066

# Row 138
$ws.Cells.Item(138, 8).Value = 55764  # H138: was 65926.336
$ws.Cells.Item(138, 10).Value = 55764  # J138: was 65926.336
$ws.Cells.Item(138, 12).Value = 55764  # L138: was 65926.336
$ws.Cells.Item(138, 14).Value = -66044  # N138: was -76206.336

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
# Row 132
$ws.Cells.Item(132, 8).Value = 1706.1351  # H132: was 1821.7059
$ws.Cells.Item(132, 9).Value = 851.75  # I132: was 929.7646999999999
$ws.Cells.Item(132, 10).Value = 2711.2942  # J132: was 2713.647
$ws.Cells.Item(132, 11).Value = 2555.25  # K132: was 2789.2941
$ws.Cells.Item(132, 12).Value = 8133.882599999999  # L132: was 8140.941
$ws.Cells.Item(132, 13).Value = -25.25  # M132: was -259.2941000000001
$ws.Cells.Item(132, 14).Value = -13193.8826  # N132: was -13200.941

# Row 136
$ws.Cells.Item(136, 8).Value = 2979268.2  # H136: was 2780744.5
$ws.Cells.Item(136, 9).Value = 3315.64  # I136: was 2989.4482
$ws.Cells.Item(136, 10).Value = 5379230  # J136: was 5379289.5
$ws.Cells.Item(136, 11).Value = 9946.92  # K136: was 8968.3446
$ws.Cells.Item(136, 12).Value = 16137690  # L136: was 16137868.5
$ws.Cells.Item(136, 13).Value = -7396.92  # M136: was -6418.3446
$ws.Cells.Item(136, 14).Value = -16142790  # N136: was -16142968.5

# Row 137
$ws.Cells.Item(137, 8).Value = 45443  # H137: was 45678.75
$ws.Cells.Item(137, 10).Value = 45443  # J137: was 45678.75
$ws.Cells.Item(137, 12).Value = 45443  # L137: was 45678.75
$ws.Cells.Item(137, 14).Value = -55643  # N137: was -55878.75

# Row 138
$ws.Cells.Item(138, 8).Value = 52114.5  # H138: was 52535.6
$ws.Cells.Item(138, 10).Value = 54229  # J138: was 54226
$ws.Cells.Item(138, 12).Value = 54229  # L138: was 54226
$ws.Cells.Item(138, 14).Value = -64509  # N138: was -64506
